# Applies updated market-board/profit figures to each job sheet's leve table.
# Generated from the authoritative OOXML diff (Sheets/Phoenix_Profits.xlsx).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 1039.091
$ws.Range("I5").Value = 1047.3334
$ws.Range("K5").Value = 1047.3334
$ws.Range("M5").Value = -932.3334

# Row 12
$ws.Range("H12").Value = 578.5
$ws.Range("I12").Value = 294.2
$ws.Range("K12").Value = 294.2
$ws.Range("M12").Value = -124.2

# Row 33
$ws.Range("H33").Value = 1565.1
$ws.Range("J33").Value = 223.5
$ws.Range("L33").Value = 223.5
$ws.Range("N33").Value = -681.5

# Row 47
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

# Row 53
$ws.Range("H53").Value = 4773.087
$ws.Range("J53").Value = 2290
$ws.Range("L53").Value = 2290
$ws.Range("N53").Value = -3564

# Row 74
$ws.Range("H74").Value = 13999.8
$ws.Range("I74").Value = 9999
$ws.Range("K74").Value = 9999
$ws.Range("M74").Value = -9063

# Row 77
$ws.Range("H77").Value = 13999.8
$ws.Range("I77").Value = 9999
$ws.Range("K77").Value = 49995
$ws.Range("M77").Value = -45315

# Row 99
$ws.Range("H99").Value = 369.54544
$ws.Range("I99").Value = 374.77777
$ws.Range("J99").Value = 346
$ws.Range("K99").Value = 1124.33331
$ws.Range("L99").Value = 1038
$ws.Range("M99").Value = 373.66669
$ws.Range("N99").Value = -4034

# Row 100
$ws.Range("H100").Value = 5781.25
$ws.Range("I100").Value = 6431.6
$ws.Range("J100").Value = 4697.3335
$ws.Range("K100").Value = 6431.6
$ws.Range("L100").Value = 4697.3335
$ws.Range("M100").Value = -5890.6
$ws.Range("N100").Value = -5779.3335

# Row 129
$ws.Range("H129").Value = 46944.24
$ws.Range("I129").Value = 91869.89999999999
$ws.Range("K129").Value = 275609.7
$ws.Range("M129").Value = -270609.7

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 847.2778
$ws.Range("I97").Value = 789.9091
$ws.Range("J97").Value = 937.4286
$ws.Range("K97").Value = 789.9091
$ws.Range("L97").Value = 937.4286
$ws.Range("M97").Value = -293.9091
$ws.Range("N97").Value = -1929.4286

# Row 128
$ws.Range("H128").Value = 183233
$ws.Range("J128").Value = 183233
$ws.Range("L128").Value = 183233
$ws.Range("N128").Value = -193193

# Row 130
$ws.Range("H130").Value = 47965.8
$ws.Range("J130").Value = 47965.8
$ws.Range("L130").Value = 47965.8
$ws.Range("N130").Value = -58005.8

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 6109.5
$ws.Range("I20").Value = 3750
$ws.Range("K20").Value = 3750
$ws.Range("M20").Value = -3503

# Row 94
$ws.Range("H94").Value = 51128.555
$ws.Range("I94").Value = 762.4545000000001
$ws.Range("J94").Value = 130275.29
$ws.Range("K94").Value = 762.4545000000001
$ws.Range("L94").Value = 130275.29
$ws.Range("M94").Value = -311.4545000000001
$ws.Range("N94").Value = -131177.29

# Row 99
$ws.Range("H99").Value = 2673.4783
$ws.Range("I99").Value = 2340
$ws.Range("J99").Value = 3298.75
$ws.Range("K99").Value = 2340
$ws.Range("L99").Value = 3298.75
$ws.Range("M99").Value = -842
$ws.Range("N99").Value = -6294.75

# Row 111
$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1793.439
$ws.Range("I31").Value = 1176.2258
$ws.Range("J31").Value = 3706.8
$ws.Range("K31").Value = 1176.2258
$ws.Range("L31").Value = 3706.8
$ws.Range("M31").Value = -881.2257999999999
$ws.Range("N31").Value = -4296.8

# Row 34
$ws.Range("H34").Value = 1793.439
$ws.Range("I34").Value = 1176.2258
$ws.Range("J34").Value = 3706.8
$ws.Range("K34").Value = 1176.2258
$ws.Range("L34").Value = 3706.8
$ws.Range("M34").Value = -974.2257999999999
$ws.Range("N34").Value = -4110.8

# Row 86
$ws.Range("H86").Value = 6403.316
$ws.Range("I86").Value = 3260.5386
$ws.Range("J86").Value = 13212.667
$ws.Range("K86").Value = 3260.5386
$ws.Range("L86").Value = 13212.667
$ws.Range("M86").Value = -2137.5386
$ws.Range("N86").Value = -15458.667

# Row 89
$ws.Range("H89").Value = 6403.316
$ws.Range("I89").Value = 3260.5386
$ws.Range("J89").Value = 13212.667
$ws.Range("K89").Value = 16302.693
$ws.Range("L89").Value = 66063.33499999999
$ws.Range("M89").Value = -10686.693
$ws.Range("N89").Value = -77295.33499999999

# Row 107
$ws.Range("H107").Value = 34918.285
$ws.Range("I107").Value = 46903.5
$ws.Range("K107").Value = 46903.5
$ws.Range("M107").Value = -44983.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 885.0714
$ws.Range("J5").Value = 907.625
$ws.Range("L5").Value = 2722.875
$ws.Range("N5").Value = -2946.875

# Row 7
$ws.Range("H7").Value = 1216.8
$ws.Range("I7").Value = 1216.8
$ws.Range("K7").Value = 3650.4
$ws.Range("M7").Value = -3538.4

# Row 12
$ws.Range("H12").Value = 829.56525
$ws.Range("I12").Value = 39.923077
$ws.Range("J12").Value = 1856.1
$ws.Range("K12").Value = 119.769231
$ws.Range("L12").Value = 5568.299999999999
$ws.Range("M12").Value = 53.23076900000001
$ws.Range("N12").Value = -5914.299999999999

# Row 23
$ws.Range("H23").Value = 222.55556
$ws.Range("I23").Value = 69
$ws.Range("K23").Value = 207
$ws.Range("M23").Value = 28

# Row 26
$ws.Range("H26").Value = 982.25
$ws.Range("I26").Value = 1340.1538
$ws.Range("J26").Value = 317.57144
$ws.Range("K26").Value = 4020.4614
$ws.Range("L26").Value = 952.71432
$ws.Range("M26").Value = -3732.4614
$ws.Range("N26").Value = -1528.71432

# Row 56
$ws.Range("H56").Value = 20000
$ws.Range("I56").Value = 20000
$ws.Range("K56").Value = 20000
$ws.Range("M56").Value = -19470

# Row 92
$ws.Range("H92").Value = 1419.75
$ws.Range("J92").Value = 1419.75
$ws.Range("L92").Value = 4259.25
$ws.Range("N92").Value = -6755.25

# Row 113
$ws.Range("H113").Value = 2446.6316
$ws.Range("J113").Value = 2642.8125
$ws.Range("L113").Value = 7928.4375
$ws.Range("N113").Value = -12268.4375

# Row 122
$ws.Range("H122").Value = 66206.57000000001
$ws.Range("I122").Value = 499.57144
$ws.Range("J122").Value = 131913.58
$ws.Range("K122").Value = 4496.14296
$ws.Range("L122").Value = 1187222.22
$ws.Range("M122").Value = -2046.14296
$ws.Range("N122").Value = -1192122.22

# Row 131
$ws.Range("H131").Value = 1871.6923
$ws.Range("I131").Value = 1004
$ws.Range("J131").Value = 3421.1428
$ws.Range("K131").Value = 3012
$ws.Range("L131").Value = 10263.4284
$ws.Range("M131").Value = 2028
$ws.Range("N131").Value = -20343.4284

# Row 132
$ws.Range("H132").Value = 2240.7778
$ws.Range("J132").Value = 4731.3335
$ws.Range("L132").Value = 42582.0015
$ws.Range("N132").Value = -47642.0015

# Row 135
$ws.Range("H135").Value = 885.0714
$ws.Range("J135").Value = 907.625
$ws.Range("L135").Value = 8168.625
$ws.Range("N135").Value = -13238.625

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7557.222
$ws.Range("I70").Value = 6005
$ws.Range("K70").Value = 6005
$ws.Range("M70").Value = -5735

# Row 73
$ws.Range("H73").Value = 7557.222
$ws.Range("I73").Value = 6005
$ws.Range("K73").Value = 6005
$ws.Range("M73").Value = -5069

# Row 97
$ws.Range("H97").Value = 1017.1667
$ws.Range("I97").Value = 732.125
$ws.Range("K97").Value = 732.125
$ws.Range("M97").Value = -236.125

# Row 101
$ws.Range("H101").Value = 19999.5
$ws.Range("J101").Value = 19999.5
$ws.Range("L101").Value = 19999.5
$ws.Range("N101").Value = -26489.5

# Row 107
$ws.Range("H107").Value = 216.47368
$ws.Range("I107").Value = 226.58333
$ws.Range("J107").Value = 199.14285
$ws.Range("K107").Value = 226.58333
$ws.Range("L107").Value = 199.14285
$ws.Range("M107").Value = 1693.41667
$ws.Range("N107").Value = -4039.14285

# Row 109
$ws.Range("H109").Value = 29988
$ws.Range("J109").Value = 29988
$ws.Range("L109").Value = 29988
$ws.Range("N109").Value = -32068

# Row 132
$ws.Range("H132").Value = 6004
$ws.Range("I132").Value = 6004
$ws.Range("K132").Value = 18012
$ws.Range("M132").Value = -15482

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3996.4285
$ws.Range("I7").Value = 3998.3333
$ws.Range("K7").Value = 3998.3333
$ws.Range("M7").Value = -3886.3333

# Row 22
$ws.Range("H22").Value = 2401.75
$ws.Range("J22").Value = 2632.6365
$ws.Range("L22").Value = 2632.6365
$ws.Range("N22").Value = -3222.6365

# Row 27
$ws.Range("H27").Value = 2401.75
$ws.Range("J27").Value = 2632.6365
$ws.Range("L27").Value = 2632.6365
$ws.Range("N27").Value = -2846.6365

# Row 40
$ws.Range("H40").Value = 6960.769
$ws.Range("I40").Value = 6203.905
$ws.Range("J40").Value = 10139.6
$ws.Range("K40").Value = 6203.905
$ws.Range("L40").Value = 10139.6
$ws.Range("M40").Value = -6067.905
$ws.Range("N40").Value = -10411.6

# Row 93
$ws.Range("H93").Value = 1798.8948
$ws.Range("I93").Value = 1344.0769
$ws.Range("J93").Value = 2784.3333
$ws.Range("K93").Value = 1344.0769
$ws.Range("L93").Value = 2784.3333
$ws.Range("M93").Value = -96.07690000000002
$ws.Range("N93").Value = -5280.3333

# Row 122
$ws.Range("H122").Value = 3219.35
$ws.Range("I122").Value = 3219.35
$ws.Range("K122").Value = 9658.049999999999
$ws.Range("M122").Value = -7208.049999999999

# Row 126
$ws.Range("H126").Value = 3996.4285
$ws.Range("I126").Value = 3998.3333
$ws.Range("K126").Value = 11994.9999
$ws.Range("M126").Value = -9524.999899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 3850
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# Row 132
$ws.Range("H132").Value = 3658.3784
$ws.Range("I132").Value = 1761.3438
$ws.Range("J132").Value = 15799.4
$ws.Range("K132").Value = 5284.0314
$ws.Range("L132").Value = 47398.2
$ws.Range("M132").Value = -2754.0314
$ws.Range("N132").Value = -52458.2
